# Link til vendor excel
# Rename the "fisk1til2"/"fisk2til3" labels to plain text "1-2"/"2-3" on the
# Ark1 (sheet1) and Ark3 (sheet3) sheets, formatting those cells as Text,
# add a new (empty, text-formatted) row below the Ark1 table, and move the
# active selection/sheet the way the author left the workbook.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Ark1")
$ws3 = $wb.Worksheets.Item("Ark3")

# --- Ark1 ("Product type" column D) -----------------------------------
$ws1.Range("D2").Value = "1-2"
$ws1.Range("D2").NumberFormat = "@"

$ws1.Range("D3").Value = "2-3"
$ws1.Range("D3").NumberFormat = "@"

$ws1.Range("D4").Value = "2-3"
$ws1.Range("D4").NumberFormat = "@"

$ws1.Range("D5").Value = "1-2"
$ws1.Range("D5").NumberFormat = "@"

# New empty row under the table, pre-formatted as text (matches the
# author adding a blank D6 cell ready for the next entry).
$ws1.Cells.Item(6, 4).NumberFormat = "@"

# --- Ark3 ("Produkttype" column B) -------------------------------------
$ws3.Range("B2").Value = "1-2"
$ws3.Range("B2").NumberFormat = "@"

$ws3.Range("B3").Value = "2-3"
$ws3.Range("B3").NumberFormat = "@"

$ws3.Range("B4").Value = "1-2"
$ws3.Range("B4").NumberFormat = "@"

$ws3.Range("B5").Value = "2-3"
$ws3.Range("B5").NumberFormat = "@"

# --- Selection / active-sheet bookkeeping ------------------------------
$ws1.Range("D6").Select()

$ws3.Activate()
$ws3.Range("E7").Select()
